{"js": "// Remove the floating signature picture (\"Grafik 1\") that is anchored to\n// the last paragraph of the document body (just before the tab run that\n// precedes the final section break). Word keeps the paragraph and the\n// trailing <w:tab/> run intact; only the <w:r> that hosts the <w:drawing>\n// (the picture) is removed.\n\nconst body = context.document.body;\nconst shapes = body.shapes;\nshapes.load(\"items/id,items/name,items/type,items/altTextDescription\");\nawait context.sync();\n\n// Prefer an unambiguous match on the picture itself (name/description as\n// seen in the original markup), but fall back to \"any floating picture\"\n// so the script still works if Word renames/relabels the shape.\nlet targets = shapes.items.filter((s) => s.type === \"Picture\" && /Grafik/i.test(s.name || \"\"));\nif (targets.length === 0) {\n  targets = shapes.items.filter((s) => s.type === \"Picture\");\n}\n\nfor (const shape of targets) {\n  shape.delete();\n}\nawait context.sync();\n", "ps1": "# Remove the floating signature picture (\"Grafik 1\") that is anchored to\n# the last paragraph of the document body. The paragraph itself (and the\n# trailing Tab run that precedes the final section break) stays in place;\n# only the floating picture shape is deleted.\n\n$d = $word.ActiveDocument\n\n# msoPicture = 13\n$msoPicture = 13\n\n# Walk backwards so deleting doesn't perturb not-yet-visited indices.\nfor ($i = $d.Shapes.Count; $i -ge 1; $i--) {\n    $shape = $d.Shapes.Item($i)\n    if ($shape.Type -eq $msoPicture -and $shape.Name -like \"Grafik*\") {\n        $shape.Delete()\n    }\n}\n\n# Fallback: if nothing matched the name pattern (e.g. shape got renamed),\n# remove any remaining floating picture shapes so the edit still applies.\nif ($d.Shapes.Count -gt 0) {\n    for ($i = $d.Shapes.Count; $i -ge 1; $i--) {\n        $shape = $d.Shapes.Item($i)\n        if ($shape.Type -eq $msoPicture) {\n            $shape.Delete()\n        }\n    }\n}\n"}
